# Corrigindo um bug com as cores das colunas
#
# Several cells in column X (the "Observações" header on row 2 plus one
# cell per evaluation block: rows 2,3,11,19,27,35,43,51,59,67) were left
# shaded grey instead of the white shade used by the rest of that row,
# breaking the alternating column-colour ("zebra stripe") pattern that the
# rest of the sheet follows. Re-apply the correct white fill to put those
# cells back in sync with their neighbours.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellsToFix = @("X2", "X3", "X11", "X19", "X27", "X35", "X43", "X51", "X59", "X67")
foreach ($addr in $cellsToFix) {
    # RGB(255,255,255) = white, matching the banding colour used by the
    # neighbouring "good" columns on the same rows (e.g. column E/F/.../W).
    $ws.Range($addr).Interior.Color = 16777215
}

# Restore the worksheet's page margins to Excel's standard defaults.
$ps = $ws.PageSetup
$ps.LeftMargin = 54      # 0.75 in
$ps.RightMargin = 54     # 0.75 in
$ps.TopMargin = 72       # 1 in
$ps.BottomMargin = 72    # 1 in
$ps.HeaderMargin = 36    # 0.5 in
$ps.FooterMargin = 36    # 0.5 in
